$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
$r1 = $tr.InsertAfter("Elkaar nakijken (schema volgt nog).")
Write-Host "r1: [$($r1.Text)] start=$($r1.Start) len=$($r1.Length)"
$nl = [char]13
$r2 = $tr.InsertAfter($nl)
Write-Host "r2: [$($r2.Text)] start=$($r2.Start) len=$($r2.Length)"
Write-Host "Full final: [$($tr.Text)]"
